$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# "TextBox 4" -- move down slightly and merge the two runs
# ("Enterprise Grid " + "Security") into a single run, text
# "Enterprise Grid Security".
$titleShape = $s.Shapes.Item(2)
$titleShape.Top = 228600 / 914400 * 72
$titleTr = $titleShape.TextFrame.TextRange
$titleTr.Delete()
$titleShape.TextFrame.TextRange.Text = "Enterprise Grid Security"

# "TextBox 9" -- drop the stray trailing endParaRPr left on the
# "http://www.cagrid.org" paragraph.
$cagridShape = $s.Shapes.Item(7)
$cagridTr = $cagridShape.TextFrame.TextRange
$cagridTr.Delete()
$cagridShape.TextFrame.TextRange.Text = "http://www.cagrid.org"

# "TextBox 10" -- drop the stray trailing endParaRPr left on the
# "GAARDS Security Infrastructure" paragraph.
$gaardsShape = $s.Shapes.Item(8)
$gaardsTr = $gaardsShape.TextFrame.TextRange
$gaardsTr.Delete()
$gaardsShape.TextFrame.TextRange.Text = "GAARDS Security Infrastructure"
